$d = $word.ActiveDocument

$d.Content.Find.Execute("648×8=5184", $true, $false, $false, $false, $false, $true, 1, $false, "723×7=5061", 2) | Out-Null
$d.Content.Find.Execute("625×9=5625", $true, $false, $false, $false, $false, $true, 1, $false, "731×5=3655", 2) | Out-Null
$d.Content.Find.Execute("365×2=730", $true, $false, $false, $false, $false, $true, 1, $false, "504×6=3024", 2) | Out-Null
$d.Content.Find.Execute("199×5=995", $true, $false, $false, $false, $false, $true, 1, $false, "656×3=1968", 2) | Out-Null
$d.Content.Find.Execute("384×5=1920", $true, $false, $false, $false, $false, $true, 1, $false, "557×4=2228", 2) | Out-Null
$d.Content.Find.Execute("919×7=6433", $true, $false, $false, $false, $false, $true, 1, $false, "661×9=5949", 2) | Out-Null
$d.Content.Find.Execute("558×8=4464", $true, $false, $false, $false, $false, $true, 1, $false, "771×8=6168", 2) | Out-Null
$d.Content.Find.Execute("286×5=1430", $true, $false, $false, $false, $false, $true, 1, $false, "294×7=2058", 2) | Out-Null
$d.Content.Find.Execute("715×3=2145", $true, $false, $false, $false, $false, $true, 1, $false, "331×7=2317", 2) | Out-Null
$d.Content.Find.Execute("572×2=1144", $true, $false, $false, $false, $false, $true, 1, $false, "143×2=286", 2) | Out-Null
$d.Content.Find.Execute("304×8=2432", $true, $false, $false, $false, $false, $true, 1, $false, "819×2=1638", 2) | Out-Null
$d.Content.Find.Execute("952×2=1904", $true, $false, $false, $false, $false, $true, 1, $false, "661×7=4627", 2) | Out-Null
$d.Content.Find.Execute("806×7=5642", $true, $false, $false, $false, $false, $true, 1, $false, "133×8=1064", 2) | Out-Null
$d.Content.Find.Execute("733×7=5131", $true, $false, $false, $false, $false, $true, 1, $false, "924×2=1848", 2) | Out-Null
$d.Content.Find.Execute("216×5=1080", $true, $false, $false, $false, $false, $true, 1, $false, "692×5=3460", 2) | Out-Null
$d.Content.Find.Execute("307×7=2149", $true, $false, $false, $false, $false, $true, 1, $false, "400×9=3600", 2) | Out-Null
$d.Content.Find.Execute("976×7=6832", $true, $false, $false, $false, $false, $true, 1, $false, "582×7=4074", 2) | Out-Null
$d.Content.Find.Execute("294×2=588", $true, $false, $false, $false, $false, $true, 1, $false, "681×6=4086", 2) | Out-Null
$d.Content.Find.Execute("228×6=1368", $true, $false, $false, $false, $false, $true, 1, $false, "500×4=2000", 2) | Out-Null
$d.Content.Find.Execute("991×8=7928", $true, $false, $false, $false, $false, $true, 1, $false, "312×3=936", 2) | Out-Null
$d.Content.Find.Execute("812×2=1624", $true, $false, $false, $false, $false, $true, 1, $false, "839×7=5873", 2) | Out-Null
$d.Content.Find.Execute("762×6=4572", $true, $false, $false, $false, $false, $true, 1, $false, "831×5=4155", 2) | Out-Null
$d.Content.Find.Execute("267×7=1869", $true, $false, $false, $false, $false, $true, 1, $false, "575×3=1725", 2) | Out-Null
$d.Content.Find.Execute("445×3=1335", $true, $false, $false, $false, $false, $true, 1, $false, "700×3=2100", 2) | Out-Null
$d.Content.Find.Execute("777×5=3885", $true, $false, $false, $false, $false, $true, 1, $false, "745×7=5215", 2) | Out-Null
